{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// The document consists of one introductory paragraph (a date line) followed\n// by a single 20x5 table of arithmetic \"answer\" cells. Each table cell holds\n// exactly one paragraph with a single run. `context.document.body.paragraphs`\n// therefore enumerates, in document order, exactly:\n//   index 0               -> the date paragraph\n//   index 1..100          -> the 100 table-cell paragraphs (row-major order)\n// which lines up 1:1 (positionally) with the values below. Several OLD\n// values repeat (e.g. \"70-9=61\" and \"82-43=39\" each appear twice) but map to\n// DIFFERENT NEW values depending on position, so replacement must be done by\n// paragraph position/index, not by global text search-and-replace.\n\nconst OLD_VALUES = [\"2024-03-14 Thursday\", \"42-39=3\", \"85-76=9\", \"49+32=81\", \"90-7=83\", \"46+38=84\", \"29+66=95\", \"19+15=34\", \"69+3=72\", \"40-9=31\", \"76+6=82\", \"39+57=96\", \"8+88=96\", \"84-68=16\", \"93-88=5\", \"15+19=34\", \"70-9=61\", \"47+47=94\", \"73-18=55\", \"48+38=86\", \"97-49=48\", \"57-8=49\", \"54+27=81\", \"29+24=53\", \"75-16=59\", \"58+13=71\", \"62-56=6\", \"50-29=21\", \"76-39=37\", \"58+37=95\", \"33-9=24\", \"56+9=65\", \"33-27=6\", \"73-68=5\", \"73-17=56\", \"62-34=28\", \"30-22=8\", \"86-27=59\", \"82-43=39\", \"57-9=48\", \"29+16=45\", \"5+78=83\", \"87+4=91\", \"28+46=74\", \"64-48=16\", \"92-68=24\", \"56-48=8\", \"92-33=59\", \"73-34=39\", \"19+38=57\", \"89+8=97\", \"93-7=86\", \"4+19=23\", \"16+45=61\", \"80-78=2\", \"64+27=91\", \"91-22=69\", \"23+8=31\", \"72-17=55\", \"51-23=28\", \"91-66=25\", \"59+37=96\", \"53+8=61\", \"7+67=74\", \"52-16=36\", \"72-64=8\", \"8+6=14\", \"6+28=34\", \"83-29=54\", \"3+48=51\", \"82-43=39\", \"22-4=18\", \"19+17=36\", \"83-58=25\", \"44-6=38\", \"26+57=83\", \"92-35=57\", \"67+14=81\", \"72-15=57\", \"7+38=45\", \"71-49=22\", \"34+29=63\", \"92-58=34\", \"19+48=67\", \"47+6=53\", \"6+9=15\", \"18+48=66\", \"90-53=37\", \"13-9=4\", \"36+57=93\", \"65+17=82\", \"29+63=92\", \"58+28=86\", \"55+8=63\", \"8+14=22\", \"64-57=7\", \"44-19=25\", \"70-9=61\", \"41-37=4\", \"95-38=57\", \"30-21=9\"];\nconst NEW_VALUES = [\"2024-03-15 Friday\", \"42+39=81\", \"36+49=85\", \"49+4=53\", \"70-38=32\", \"59+33=92\", \"40-12=28\", \"17+54=71\", \"16+78=94\", \"86-67=19\", \"34-27=7\", \"3+78=81\", \"29+59=88\", \"73-14=59\", \"78-39=39\", \"66+18=84\", \"19+68=87\", \"71-16=55\", \"33-17=16\", \"15+36=51\", \"90-1=89\", \"29+22=51\", \"69+29=98\", \"9+75=84\", \"77+5=82\", \"24-19=5\", \"83-17=66\", \"97-78=19\", \"91-65=26\", \"29+62=91\", \"81-6=75\", \"32-25=7\", \"84-48=36\", \"95-57=38\", \"77-49=28\", \"65+8=73\", \"52-6=46\", \"66-57=9\", \"74-55=19\", \"12-4=8\", \"17+39=56\", \"95-9=86\", \"67+25=92\", \"80-31=49\", \"45+19=64\", \"77-59=18\", \"6+55=61\", \"35-29=6\", \"36+49=85\", \"9+58=67\", \"48+18=66\", \"50-45=5\", \"94-27=67\", \"43-14=29\", \"93-64=29\", \"19+54=73\", \"88-19=69\", \"70-31=39\", \"53-38=15\", \"73-65=8\", \"46+6=52\", \"88+5=93\", \"91-79=12\", \"57-48=9\", \"91-16=75\", \"80-3=77\", \"81-78=3\", \"18+68=86\", \"25+37=62\", \"60-32=28\", \"46+26=72\", \"9+13=22\", \"68+29=97\", \"33-16=17\", \"83-67=16\", \"52-46=6\", \"50-38=12\", \"17+15=32\", \"90-14=76\", \"90-88=2\", \"78+18=96\", \"38+28=66\", \"12+59=71\", \"19+75=94\", \"69+13=82\", \"2+79=81\", \"39+35=74\", \"39+16=55\", \"14+48=62\", \"37-8=29\", \"27+59=86\", \"52-43=9\", \"61-52=9\", \"8+86=94\", \"72-13=59\", \"9+18=27\", \"45+8=53\", \"90-35=55\", \"44-29=15\", \"92-34=58\", \"70-18=52\"];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nif (paragraphs.items.length !== NEW_VALUES.length) {\n  throw new Error(\n    \"Unexpected paragraph count: expected \" + NEW_VALUES.length +\n    \" but found \" + paragraphs.items.length\n  );\n}\n\n// Load current text for each paragraph so we can sanity-check before writing.\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const expectedOld = OLD_VALUES[i];\n  const newValue = NEW_VALUES[i];\n  const currentText = para.text;\n\n  if (currentText === newValue) {\n    // Already at the target value; nothing to do.\n    continue;\n  }\n\n  if (currentText !== expectedOld) {\n    console.log(\n      \"WARNING: paragraph \" + i + \" text was \" + JSON.stringify(currentText) +\n      \", expected \" + JSON.stringify(expectedOld)\n    );\n  }\n\n  // `Paragraph.text` is read-only in Office.js; replace the paragraph's\n  // content via insertText(..., \"Replace\") to preserve run formatting\n  // (font, size, etc.) on the existing run.\n  para.insertText(newValue, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# PowerShell / Word COM-interop edit script.\n# Operates on $word.ActiveDocument ($d below).\n#\n# Layout: one intro paragraph (a date line) followed by a single 20x5 table\n# of arithmetic \"answer\" cells (one paragraph/run per cell). Several OLD cell\n# values repeat (e.g. \"70-9=61\", \"82-43=39\" each appear twice) but map to\n# DIFFERENT NEW values depending on position, so every cell is addressed\n# explicitly by (row, column) rather than via a global Find/Replace.\n\n$d = $word.ActiveDocument\n\n# --- Date paragraph -------------------------------------------------------\n$expectedDateOld = '2024-03-14 Thursday'\n$dateNew = '2024-03-15 Friday'\n$p1 = $d.Paragraphs(1)\n$p1Text = $p1.Range.Text.TrimEnd([char]13, [char]7)\nif ($p1Text -ne $expectedDateOld) {\n    Write-Output \"WARNING: date paragraph text was '$p1Text', expected '$expectedDateOld'\"\n}\n$p1.Range.Text = $dateNew\n\n# --- Table cells (20 rows x 5 columns), row-major, 1-based indices --------\n$t = $d.Tables(1)\n\n$cellValues = @(\n    @(@('42-39=3','42+39=81'), @('85-76=9','36+49=85'), @('49+32=81','49+4=53'), @('90-7=83','70-38=32'), @('46+38=84','59+33=92')),\n    @(@('29+66=95','40-12=28'), @('19+15=34','17+54=71'), @('69+3=72','16+78=94'), @('40-9=31','86-67=19'), @('76+6=82','34-27=7')),\n    @(@('39+57=96','3+78=81'), @('8+88=96','29+59=88'), @('84-68=16','73-14=59'), @('93-88=5','78-39=39'), @('15+19=34','66+18=84')),\n    @(@('70-9=61','19+68=87'), @('47+47=94','71-16=55'), @('73-18=55','33-17=16'), @('48+38=86','15+36=51'), @('97-49=48','90-1=89')),\n    @(@('57-8=49','29+22=51'), @('54+27=81','69+29=98'), @('29+24=53','9+75=84'), @('75-16=59','77+5=82'), @('58+13=71','24-19=5')),\n    @(@('62-56=6','83-17=66'), @('50-29=21','97-78=19'), @('76-39=37','91-65=26'), @('58+37=95','29+62=91'), @('33-9=24','81-6=75')),\n    @(@('56+9=65','32-25=7'), @('33-27=6','84-48=36'), @('73-68=5','95-57=38'), @('73-17=56','77-49=28'), @('62-34=28','65+8=73')),\n    @(@('30-22=8','52-6=46'), @('86-27=59','66-57=9'), @('82-43=39','74-55=19'), @('57-9=48','12-4=8'), @('29+16=45','17+39=56')),\n    @(@('5+78=83','95-9=86'), @('87+4=91','67+25=92'), @('28+46=74','80-31=49'), @('64-48=16','45+19=64'), @('92-68=24','77-59=18')),\n    @(@('56-48=8','6+55=61'), @('92-33=59','35-29=6'), @('73-34=39','36+49=85'), @('19+38=57','9+58=67'), @('89+8=97','48+18=66')),\n    @(@('93-7=86','50-45=5'), @('4+19=23','94-27=67'), @('16+45=61','43-14=29'), @('80-78=2','93-64=29'), @('64+27=91','19+54=73')),\n    @(@('91-22=69','88-19=69'), @('23+8=31','70-31=39'), @('72-17=55','53-38=15'), @('51-23=28','73-65=8'), @('91-66=25','46+6=52')),\n    @(@('59+37=96','88+5=93'), @('53+8=61','91-79=12'), @('7+67=74','57-48=9'), @('52-16=36','91-16=75'), @('72-64=8','80-3=77')),\n    @(@('8+6=14','81-78=3'), @('6+28=34','18+68=86'), @('83-29=54','25+37=62'), @('3+48=51','60-32=28'), @('82-43=39','46+26=72')),\n    @(@('22-4=18','9+13=22'), @('19+17=36','68+29=97'), @('83-58=25','33-16=17'), @('44-6=38','83-67=16'), @('26+57=83','52-46=6')),\n    @(@('92-35=57','50-38=12'), @('67+14=81','17+15=32'), @('72-15=57','90-14=76'), @('7+38=45','90-88=2'), @('71-49=22','78+18=96')),\n    @(@('34+29=63','38+28=66'), @('92-58=34','12+59=71'), @('19+48=67','19+75=94'), @('47+6=53','69+13=82'), @('6+9=15','2+79=81')),\n    @(@('18+48=66','39+35=74'), @('90-53=37','39+16=55'), @('13-9=4','14+48=62'), @('36+57=93','37-8=29'), @('65+17=82','27+59=86')),\n    @(@('29+63=92','52-43=9'), @('58+28=86','61-52=9'), @('55+8=63','8+86=94'), @('8+14=22','72-13=59'), @('64-57=7','9+18=27')),\n    @(@('44-19=25','45+8=53'), @('70-9=61','90-35=55'), @('41-37=4','44-29=15'), @('95-38=57','92-34=58'), @('30-21=9','70-18=52'))\n)\n\nfor ($r = 1; $r -le 20; $r++) {\n    for ($c = 1; $c -le 5; $c++) {\n        $pair = $cellValues[$r - 1][$c - 1]\n        $expectedOld = $pair[0]\n        $newVal = $pair[1]\n        $cell = $t.Cell($r, $c)\n        $cellText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n        if ($cellText -ne $expectedOld) {\n            Write-Output \"WARNING: cell ($r,$c) text was '$cellText', expected '$expectedOld'\"\n        }\n        $cell.Range.Text = $newVal\n    }\n}\n\nWrite-Output 'done'\n"}
